# "Generate Report for Archive"
#
# The two tracked files (721aa51a-...md / af6c4662-...md) swapped places in
# the localization-status report: the row that used to describe
# 721aa51a-...md now describes af6c4662-...md (and vice versa), and the
# "Status" column for that pair moved from "Ready for handoff" to
# "In Translation". This touches the Overview sheet plus the per-locale
# (zh-cn / de-de) detail sheets, including the hyperlinks that decorate the
# file-name columns.

$wb = $excel.ActiveWorkbook

$id721 = "721aa51a-e9be-4dc0-9833-32873f099577.md"
$idaf6 = "af6c4662-f8fd-4e34-957c-3654765d9d23.md"

$url721md = "https://github.com/OpenLocalizationTest/oltest/blob/eed213af7ca40fd417abe975e67be74d227528f8/e2e/721aa51a-e9be-4dc0-9833-32873f099577.md"
$urlaf6md = "https://github.com/OpenLocalizationTest/oltest/blob/eed213af7ca40fd417abe975e67be74d227528f8/e2e/af6c4662-f8fd-4e34-957c-3654765d9d23.md"
$urlConfig = "https://github.com/OpenLocalizationTest/oltest/blob/eed213af7ca40fd417abe975e67be74d227528f8/.localization-config"

$newStatus = "In Translation"

# ---------------------------------------------------------------------
# Overview sheet: A2/A3 file names swap, B/C status becomes "In Translation"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $idaf6
$wsOverview.Range("A3").Value = $id721
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $url721md, "", "", $idaf6)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $urlaf6md, "", "", $id721)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $urlConfig, "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$urlZhXlf721 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e384b61201ce780c9dd60048116ca64bb0b41c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/721aa51a-e9be-4dc0-9833-32873f099577.3fba059ee59d5fa5ed5fd5aa2effe57f558ed525.zh-cn.xlf"
$urlZhXlfaf6 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e384b61201ce780c9dd60048116ca64bb0b41c7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.zh-cn.xlf"

$zhXlf721 = "721aa51a-e9be-4dc0-9833-32873f099577.3fba059ee59d5fa5ed5fd5aa2effe57f558ed525.zh-cn.xlf"
$zhXlfaf6 = "af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.zh-cn.xlf"

$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $idaf6
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Value = $zhXlfaf6
$wsZh.Range("A3").Value = $id721
$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("C3").Value = $zhXlf721

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $url721md, "", "", $idaf6)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $urlZhXlf721, "", "", $zhXlfaf6)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlaf6md, "", "", $id721)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $urlZhXlfaf6, "", "", $zhXlf721)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $urlConfig, "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$urlDeXlf721 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32713c9abb62d7025c31384c79b02b15274b5191/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/721aa51a-e9be-4dc0-9833-32873f099577.3fba059ee59d5fa5ed5fd5aa2effe57f558ed525.de-de.xlf"
$urlDeXlfaf6 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32713c9abb62d7025c31384c79b02b15274b5191/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.de-de.xlf"

$deXlf721 = "721aa51a-e9be-4dc0-9833-32873f099577.3fba059ee59d5fa5ed5fd5aa2effe57f558ed525.de-de.xlf"
$deXlfaf6 = "af6c4662-f8fd-4e34-957c-3654765d9d23.3759a12535d2c9f4036116f9969abb4278de4a85.de-de.xlf"

$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $idaf6
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Value = $deXlfaf6
$wsDe.Range("A3").Value = $id721
$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("C3").Value = $deXlf721

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $url721md, "", "", $idaf6)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $urlDeXlf721, "", "", $deXlfaf6)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlaf6md, "", "", $id721)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $urlDeXlfaf6, "", "", $deXlf721)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $urlConfig, "", "", ".localization-config")
